# Fruta / hortaliza, semanal
# A new weekly price observation was added for "Agrícola del Norte S.A. de
# Arica - Kiwi": insert a new row above the current row 13, pushing the
# existing rows 13-15 down to 14-16, and populate the new row 13 with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13:15 shift down to 14:16,
# carrying their values/formatting with them (matches the diff, where the
# old row 13 data reappears verbatim as the new row 14, etc.).
$ws.Rows("13:13").Insert()

# Fill in the newly inserted row 13 with this week's observation.
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = 44706
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100101
$ws.Range("H13").Value = "Berries"
$ws.Range("I13").Value = 100101007
$ws.Range("J13").Value = "Kiwi"
$ws.Range("K13").Value = "Hayward"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 400
$ws.Range("N13").Value = 9000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 9500
$ws.Range("Q13").Value = "$/bandeja 10 kilos"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 950
$ws.Range("T13").Value = 10
